# Change layout of utilization tab:
# Clear contents of A7:C11 (leave formatting/styles intact), and update
# the active cell selection to I12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7:C11").ClearContents()

$ws.Range("I12").Select()
